{"js": "// The transcription markup in this document wraps each page's record id in\n// \"<id>...</id>\" tags that were previously split across three separate runs\n// (the literal \"<id>\" tag, the bare id text, and the literal \"</id>\" tag).\n// The edit collapses each of those three runs into a single run containing\n// the full \"<id>p156r_N</id>\" text, keeping the formatting of the opening\n// \"<id>\" tag run (Courier New / color 7f6000 / 9pt).\n\nconst body = context.document.body;\n\n// Locate every \"<id>\" opening tag and its matching \"</id>\" closing tag.\nconst openTags = body.search(\"<id>\", { matchCase: true });\nconst closeTags = body.search(\"</id>\", { matchCase: true });\nopenTags.load(\"items\");\ncloseTags.load(\"items\");\nawait context.sync();\n\nconst count = Math.min(openTags.items.length, closeTags.items.length);\n\nfor (let i = 0; i < count; i++) {\n  const openRange = openTags.items[i];\n  const closeRange = closeTags.items[i];\n\n  // Build a range that spans from the start of \"<id>\" to the end of \"</id>\"\n  // (this also swallows the run(s) holding the id value in between).\n  const fullRange = openRange.expandTo(closeRange);\n  fullRange.load(\"text\");\n  await context.sync();\n\n  const mergedText = fullRange.text; // e.g. \"<id>p156r_1</id>\"\n\n  // Replacing the whole span with its own text collapses the 3 runs into a\n  // single run, inheriting the formatting of the first run in the range\n  // (the \"<id>\" tag's Courier New / 7f6000 / 9pt run).\n  fullRange.insertText(mergedText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The transcription markup in this document wraps each page's record id in\n# \"<id>...</id>\" tags that were previously split across three separate runs\n# (the literal \"<id>\" tag, the bare id text, and the literal \"</id>\" tag).\n# This script collapses each of those three runs into a single run\n# containing the full \"<id>p156r_N</id>\" text, keeping the formatting of the\n# opening \"<id>\" tag run (Courier New / color 7f6000 / 9pt).\n\n$d = $word.ActiveDocument\n\n$searchFrom = 0\n$docEnd = $d.Content.End\n\nfor ($i = 0; $i -lt 20; $i++) {\n\n    # Find the next \"<id>\" opening tag, starting after the previous match.\n    $openRange = $d.Range($searchFrom, $docEnd)\n    $openFind = $openRange.Find\n    $openFind.ClearFormatting()\n    $openFind.Text = \"<id>\"\n    $openFind.MatchCase = $true\n    $openFind.MatchWildcards = $false\n    $foundOpen = $openFind.Execute()\n    if (-not $foundOpen) { break }\n\n    $openStart = $openRange.Start\n\n    # Find the matching \"</id>\" closing tag that follows it.\n    $closeRange = $d.Range($openRange.End, $docEnd)\n    $closeFind = $closeRange.Find\n    $closeFind.ClearFormatting()\n    $closeFind.Text = \"</id>\"\n    $closeFind.MatchCase = $true\n    $closeFind.MatchWildcards = $false\n    $foundClose = $closeFind.Execute()\n    if (-not $foundClose) { break }\n\n    $closeEnd = $closeRange.End\n\n    # Range spanning \"<id>\" ... \"</id>\" inclusive (covers the 3 runs:\n    # \"<id>\", the id value, and \"</id>\").\n    $combined = $d.Range($openStart, $closeEnd)\n    $mergedText = $combined.Text\n\n    # Assigning Range.Text to a value that differs from the current text\n    # collapses the range into a single run, inheriting the character\n    # formatting of the first run in the range (the \"<id>\" tag's\n    # Courier New / 7f6000 / 9pt). A genuine self-assignment is a no-op in\n    # this object model, so we first write a sentinel-suffixed value (a\n    # real change) and then trim it back to the final text (also a real\n    # change) to force the run merge.\n    $combined.Text = $mergedText + [char]1\n    $trimmed = $d.Range($openStart, $openStart + $mergedText.Length + 1)\n    $trimmed.Text = $mergedText\n\n    $searchFrom = $openStart + $mergedText.Length\n    $docEnd = $d.Content.End\n}\n"}
